# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1)
# onto the three new header cells, then set their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-47) gets the team's season record: 95 wins, 67
# losses, 0 ties.
$ws.Range("AD2:AD47").Value = 95
$ws.Range("AE2:AE47").Value = 67
$ws.Range("AF2:AF47").Value = 0

$excel.CutCopyMode = $false
